$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Handout master: the auto date field cached text "4/22/20" -> "5/7/20"
#    (PowerPoint re-stamps the datetimeFigureOut field's cached text on
#    every save). The handout master's date placeholder is reachable via
#    HeadersFooters.DateAndTime on the HandoutMaster.
# ---------------------------------------------------------------------
$hm = $p.HandoutMaster
$dt = $hm.HeadersFooters.DateAndTime
$dt.Text = "5/7/20"

# ---------------------------------------------------------------------
# 2) Slide 3 ("Requirements and Scope"): bullet line rewording.
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3Body = $s3.Shapes.Item(2).TextFrame.TextRange
$old3 = "RFC 5357 (TWAMP) defined probe messages - TWAMP Light"
$new3 = "RFC 5357 (TWAMP Light) defined probe messages"
$i3 = $s3Body.Text.IndexOf($old3)
$s3Body.Characters($i3 + 1, $old3.Length).Text = $new3

# ---------------------------------------------------------------------
# 3) Slide 4 ("History of the Draft"): un-bold the "SPRING Chairs..." run.
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$s4Body = $s4.Shapes.Item(2).TextFrame.TextRange
$bold4 = "SPRING Chairs announced in the meeting the agreement with IPPM chairs to progress the draft in SPRING WG"
$i4 = $s4Body.Text.IndexOf($bold4)
$s4Body.Characters($i4 + 1, $bold4.Length).Font.Bold = 0

# ---------------------------------------------------------------------
# 4) Slide 12 ("Probe Query for Links"): reword the DM/LM probe sentence,
#    keeping the bold "port1"/"port2" runs untouched. Replace the later
#    substring first so the earlier substring's offset stays valid.
# ---------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$s12Body = $s12.Shapes.Item(3).TextFrame.TextRange

$oldLM = " is used for LM probe messages in unauthenticated mode."
$newLM = " is used for LM probe messages, both in unauthenticated mode."
$iLM = $s12Body.Text.IndexOf($oldLM)
$s12Body.Characters($iLM + 1, $oldLM.Length).Text = $newLM

$oldDM = " is used for DM probe messages in unauthenticated mode and "
$newDM = " is used for DM probe messages and "
$iDM = $s12Body.Text.IndexOf($oldDM)
$s12Body.Characters($iDM + 1, $oldDM.Length).Text = $newDM
